$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Punters 2025": update round results for R16-R19 (rows 20-23).
# Downstream totals (row 29), the summary block (rows 31-38) and the
# cross-sheet comparison block (AI7:AL18, which reads from this sheet's
# column totals) are all formula-driven, so they recalculate automatically.
# ---------------------------------------------------------------------------
$p25 = $wb.Worksheets.Item("Punters 2025")

# Row 20 (R16)
$p25.Range("K20").Value = 7.5
$p25.Range("L20").Value = 0
$p25.Range("O20").Value = 5
$p25.Range("P20").Value = 0
$p25.Range("S20").Value = 25
$p25.Range("T20").Value = 0
$p25.Range("U20").Value = 25
$p25.Range("V20").Value = 0
$p25.Range("W20").Value = 25
$p25.Range("X20").Value = 0
$p25.Range("Y20").Value = 25
$p25.Range("Z20").Value = 0

# Row 21 (R17)
$p25.Range("C21").Formula = "=6+5+4.5+4.5"
$p25.Range("D21").Value = 45
$p25.Range("E21").Value = 15
$p25.Range("F21").Value = 0
$p25.Range("G21").Formula = "=10+1+5+10"
$p25.Range("H21").Value = 0
$p25.Range("I21").Formula = "=5+5+5"
$p25.Range("J21").Value = 0

# Row 22 (R18)
$p25.Range("C22").Value = 5
$p25.Range("D22").Value = 0
$p25.Range("I22").Value = 10
$p25.Range("J22").Value = 0
$p25.Range("K22").Value = 17.5
$p25.Range("L22").Value = 0
$p25.Range("O22").Value = 25
$p25.Range("P22").Value = 0
$p25.Range("Q22").Value = 25
$p25.Range("R22").Value = 0

# Row 23 (R19)
$p25.Range("K23").Value = 7.5
$p25.Range("L23").Value = 0
$p25.Range("M23").Value = 20
$p25.Range("N23").Value = 0
$p25.Range("S23").Value = 25
$p25.Range("T23").Value = 0
$p25.Range("U23").Value = 25
$p25.Range("V23").Value = 0
$p25.Range("W23").Value = 25
$p25.Range("X23").Value = 69.83
$p25.Range("Y23").Value = 27
$p25.Range("Z23").Value = 0

# ---------------------------------------------------------------------------
# Sheet "1 Leg Losses 2025": append this week's new 1-leg-loss bets (rows
# 38-44). Column F (per-punter totals, rows 1-13) is formula-driven and
# recalculates automatically.
# ---------------------------------------------------------------------------
$losses = $wb.Worksheets.Item("1 Leg Losses 2025")

$dollarCents  = """$""#,##0.00_);[Red]\(""$""#,##0.00\)"
$dollarWhole  = """$""#,##0_);[Red]\(""$""#,##0\)"

$losses.Range("A38").Value = "Fake"
$losses.Range("B38").Value = "`$5 bonus"
$losses.Range("C38").Value = 95
$losses.Range("C38").NumberFormat = $dollarWhole
$losses.Range("D38").Value = "Max gawn 14/15 touches"

$losses.Range("A39").Value = "BT"
$losses.Range("B39").Value = 7.5
$losses.Range("B39").NumberFormat = $dollarCents
$losses.Range("C39").Value = 63.75
$losses.Range("C39").NumberFormat = $dollarCents
$losses.Range("D39").Value = "Max gawn 14/15 touches"

$losses.Range("A40").Value = "Simmo"
$losses.Range("B40").Value = 7
$losses.Range("B40").NumberFormat = $dollarWhole
$losses.Range("C40").Value = 106.75
$losses.Range("C40").NumberFormat = $dollarCents
$losses.Range("D40").Value = "Caleb Serong most disposals"

$losses.Range("A41").Value = "James"
$losses.Range("B41").Value = 15
$losses.Range("B41").NumberFormat = $dollarWhole
$losses.Range("C41").Value = 157.5
$losses.Range("C41").NumberFormat = $dollarCents
$losses.Range("D41").Value = "Petracca 23/24 disposals"

$losses.Range("A42").Value = "Scott"
$losses.Range("B42").Value = 10
$losses.Range("B42").NumberFormat = $dollarWhole
$losses.Range("C42").Value = 105
$losses.Range("C42").NumberFormat = $dollarWhole
$losses.Range("D42").Value = "MucCluggage 18/25.5 disposals"

$losses.Range("A43").Value = "Fake"
$losses.Range("B43").Value = 5
$losses.Range("B43").NumberFormat = $dollarWhole
$losses.Range("C43").Value = 166.25
$losses.Range("C43").NumberFormat = $dollarCents
$losses.Range("D43").Value = "Herbie Farnworth anytime try scorer"

$losses.Range("A44").Value = "BT"
$losses.Range("B44").Value = 7.5
$losses.Range("B44").NumberFormat = $dollarCents
$losses.Range("C44").Value = 90
$losses.Range("C44").NumberFormat = $dollarWhole
$losses.Range("D44").Value = "GWS Essendon over 171.5, total 160"

# ---------------------------------------------------------------------------
# Text corrections in shared strings (typo fix + date updates).
# ---------------------------------------------------------------------------
$wb.Worksheets | ForEach-Object {
    $used = $_.UsedRange
    $found = $used.Find("Swikowski 6/15 touches")
    if ($found -ne $null) { $found.Value = "Switkowski 6/15 touches" }

    $found2 = $used.Find("1 leg losses as at 16 June 2025")
    if ($found2 -ne $null) { $found2.Value = "1 leg losses as at 23 July 2025" }

    $found3 = $used.Find("Strikes for 2025 Season at 10 June 2025")
    if ($found3 -ne $null) { $found3.Value = "Strikes for 2025 Season at 24 July 2025" }
}
